$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Snca"
$ws.Range("C2").Value = "Lag3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04830066666666667
$ws.Range("H2").Value = 0.144902
$ws.Range("I2").Value = 0.04634491140536046
$ws.Range("J2").Value = 0.04634491140536046
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.574702666666667
$ws.Range("N2").Value = 22.724108
$ws.Range("O2").Value = 0.2360813295275979
$ws.Range("P2").Value = 0.2360813295275979
$ws.Range("Q2").Value = 0.3658631886017778
$ws.Range("R2").Value = 3.292768697416
$ws.Range("S2").Value = 0.01094116830141623
$ws.Range("T2").Value = 0.01094116830141623

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Snca"
$ws.Range("C3").Value = "Lag3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04830066666666667
$ws.Range("H3").Value = 0.144902
$ws.Range("I3").Value = 0.04634491140536046
$ws.Range("J3").Value = 0.04634491140536046
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.86426066666667
$ws.Range("N3").Value = 56.592782
$ws.Range("O3").Value = 0.5879438355171306
$ws.Range("P3").Value = 0.5879438355171307
$ws.Range("Q3").Value = 0.9111563663737777
$ws.Range("R3").Value = 8.200407297364
$ws.Range("S3").Value = 0.02724820496836924
$ws.Range("T3").Value = 0.02724820496836924

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Snca"
$ws.Range("C4").Value = "Lag3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04830066666666667
$ws.Range("H4").Value = 0.144902
$ws.Range("I4").Value = 0.04634491140536046
$ws.Range("J4").Value = 0.04634491140536046
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.646177333333333
$ws.Range("N4").Value = 16.938532
$ws.Range("O4").Value = 0.1759748349552714
$ws.Range("P4").Value = 0.1759748349552714
$ws.Range("Q4").Value = 0.2727141293182222
$ws.Range("R4").Value = 2.454427163864
$ws.Range("S4").Value = 0.00815553813557498
$ws.Range("T4").Value = 0.00815553813557498

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Snca"
$ws.Range("C5").Value = "Lag3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9938993333333332
$ws.Range("H5").Value = 2.981698
$ws.Range("I5").Value = 0.9536550885946394
$ws.Range("J5").Value = 0.9536550885946395
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.574702666666667
$ws.Range("N5").Value = 22.724108
$ws.Range("O5").Value = 0.2360813295275979
$ws.Range("P5").Value = 0.2360813295275979
$ws.Range("Q5").Value = 7.528491930598222
$ws.Range("R5").Value = 67.75642737538399
$ws.Range("S5").Value = 0.2251401612261817
$ws.Range("T5").Value = 0.2251401612261817

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Snca"
$ws.Range("C6").Value = "Lag3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9938993333333332
$ws.Range("H6").Value = 2.981698
$ws.Range("I6").Value = 0.9536550885946394
$ws.Range("J6").Value = 0.9536550885946395
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.86426066666667
$ws.Range("N6").Value = 56.592782
$ws.Range("O6").Value = 0.5879438355171306
$ws.Range("P6").Value = 0.5879438355171307
$ws.Range("Q6").Value = 18.74917610042622
$ws.Range("R6").Value = 168.742584903836
$ws.Range("S6").Value = 0.5606956305487613
$ws.Range("T6").Value = 0.5606956305487615

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Snca"
$ws.Range("C7").Value = "Lag3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9938993333333332
$ws.Range("H7").Value = 2.981698
$ws.Range("I7").Value = 0.9536550885946394
$ws.Range("J7").Value = 0.9536550885946395
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.646177333333333
$ws.Range("N7").Value = 16.938532
$ws.Range("O7").Value = 0.1759748349552714
$ws.Range("P7").Value = 0.1759748349552714
$ws.Range("Q7").Value = 5.611731887481777
$ws.Range("R7").Value = 50.50558698733599
$ws.Range("S7").Value = 0.1678192968196964
$ws.Range("T7").Value = 0.1678192968196964
